$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new "Instalacion" ticket in TEMUCO (Panamericana Sur) ---
$ws.Range("A2").Value = 46021
$ws.Range("B2").Value = 43572
$ws.Range("C2").Value = "normal"
$ws.Range("D2").Value = "Instalación"
$ws.Range("E2").Value = "GPS, Corta Corriente, Sensor Pta, Sensor Temperatura"
$ws.Range("F2").Value = "PANAMERICANA SUR KM. 678 PADRE LAS CASAS"
$ws.Range("G2").Value = "TEMUCO"
$ws.Range("H2").Value = "Región de La Araucanía."
$ws.Range("I2").Value = "Pedro Prez"
$ws.Range("J2").Value = "TRPT29"
$ws.Range("K2").Value = "DIWATTS"
$ws.Range("L2").Value = "GPS"
$ws.Range("M2").Value = " Corta Corriente"
$ws.Range("N2").Value = " Sensor Pta"
$ws.Range("O2").Value = " Sensor Temperatura"
$ws.Range("P2:R2").ClearContents()

# --- Row 3: "Instalacion" ticket in CHILLAN (Ruta 5 Sur) ---
$ws.Range("A3").Value = 46021
$ws.Range("B3").Value = 42426
$ws.Range("C3").Value = "normal"
$ws.Range("D3").Value = "Instalación"
$ws.Range("E3").Value = "GPS, Botón Alámbrico Tablero, Corta Corriente, Sensor Pta, Sensor Pta Adicional, Sensor Temperatura, Sensor Temperatura Adicional"
$ws.Range("F3").Value = "RUTA 5 SUR KM 8 SN"
$ws.Range("G3").Value = "CHILLAN"
$ws.Range("H3").Value = "Región del Ñuble."
$ws.Range("I3").Value = "Pedro Pascal"
$ws.Range("J3").Value = "VE839-POR CONFIRMAR"
$ws.Range("K3").Value = "CIAL_ALIMENTOS"
$ws.Range("L3").Value = "GPS"
$ws.Range("M3").Value = " Botón Alámbrico Tablero"
$ws.Range("N3").Value = " Corta Corriente"
$ws.Range("O3").Value = " Sensor Pta"
$ws.Range("P3").Value = " Sensor Pta Adicional"
$ws.Range("Q3").Value = " Sensor Temperatura"
$ws.Range("R3").Value = " Sensor Temperatura Adicional"

# --- Row 4 (new): "Soporte" ticket in SAN BERNARDO ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 46021
$ws.Range("B4").Value = 43567
$ws.Range("C4").Value = "normal"
$ws.Range("D4").Value = "Soporte"
$ws.Range("E4").Value = "GPS"
$ws.Range("F4").Value = "CAMINO LONGITUDINAL SUR 5201, NOS"
$ws.Range("G4").Value = "SAN BERNARDO"
$ws.Range("H4").Value = "Región Metropolitana de Santiago."
$ws.Range("I4").Value = "Pedro Prez"
$ws.Range("J4").Value = "HLPX63"
$ws.Range("K4").Value = "CUENTA_CAROZZIDISTRIBUCION"
$ws.Range("L4").Value = "GPS"

# --- Row 5 (new): "Soporte" ticket in TEMUCO ---
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 46021
$ws.Range("B5").Value = 43094
$ws.Range("C5").Value = "normal"
$ws.Range("D5").Value = "Soporte"
$ws.Range("E5").Value = "GPS"
$ws.Range("F5").Value = "GUIDO BECK DE RAMBERGA 1884, PADRE DE LAS CASAS "
$ws.Range("G5").Value = "TEMUCO"
$ws.Range("H5").Value = "Región de La Araucanía."
$ws.Range("I5").Value = "Pedro Pascal"
$ws.Range("J5").Value = "BWYY79"
$ws.Range("K5").Value = "CUENTA_CAROZZIDISTRIBUCION"
$ws.Range("L5").Value = "GPS"

$ws.Range("C4").Select()
